# Generate Report for Handback
# Updates timestamp / priority values produced by a fresh handback report run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# "Latest HO Xliff Generate Date" for 72e53071-... and 7ca14f11-... rows on Overview
$wsOverview.Range("G3").Value = "2016-08-20 14:13:15"
$wsOverview.Range("G4").Value = "2016-08-20 14:13:15"

# "Priority" changed from ht to mt for the same two rows, in both locale sheets
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"

# "Correspond Handoff Datetime" (zh-cn) for the same two rows
$wsZhCn.Range("H3").Value = "2016-08-20 14:13:11"
$wsZhCn.Range("H4").Value = "2016-08-20 14:13:11"

# "Correspond Handback DateTime" (zh-cn) for the same two rows
$wsZhCn.Range("K3").Value = "2016-08-20 14:13:27"
$wsZhCn.Range("K4").Value = "2016-08-20 14:13:27"

# "Correspond Handoff Datetime" (de-de) for the same two rows (mirrors Overview G3/G4)
$wsDeDe.Range("H3").Value = "2016-08-20 14:13:15"
$wsDeDe.Range("H4").Value = "2016-08-20 14:13:15"

# "Correspond Handback DateTime" (de-de) for the same two rows
$wsDeDe.Range("K3").Value = "2016-08-20 14:13:33"
$wsDeDe.Range("K4").Value = "2016-08-20 14:13:33"
